$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 2553
$ws.Range("E4").Value = 11113
$ws.Range("D5").Value = "1066 (41.8)"
$ws.Range("E5").Value = "4222 (38.0)"
$ws.Range("D6").Value = "623 (24.4)"
$ws.Range("E6").Value = "2677 (24.1)"
$ws.Range("D7").Value = "526 (20.6)"
$ws.Range("E7").Value = "2496 (22.5)"
$ws.Range("D8").Value = "338 (13.2)"
$ws.Range("E8").Value = "1718 (15.5)"
$ws.Range("D9").Value = "1217 (47.7)"
$ws.Range("E9").Value = "4665 (42.0)"
$ws.Range("D10").Value = "973 (38.1)"
$ws.Range("E10").Value = "5641 (50.8)"
$ws.Range("D11").Value = "357 (14.0)"
$ws.Range("E11").Value = "594 (5.3)"
$ws.Range("D12").Value = "1223 (47.9)"
$ws.Range("E12").Value = "4878 (43.9)"
$ws.Range("D13").Value = "803 (31.5)"
$ws.Range("E13").Value = "518 (4.7)"
$ws.Range("D14").Value = "200 (7.8)"
$ws.Range("E14").Value = "1576 (14.2)"
$ws.Range("D15").Value = "835 (32.7)"
$ws.Range("E15").Value = "4639 (41.7)"
$ws.Range("D16").Value = "503 (19.7)"
$ws.Range("E16").Value = "2082 (18.7)"
$ws.Range("D17").Value = "726 (28.4)"
$ws.Range("E17").Value = "3055 (27.5)"
$ws.Range("D18").Value = "1282 (50.2)"
$ws.Range("E18").Value = "5712 (51.4)"
$ws.Range("D19").Value = "244 (9.6)"
$ws.Range("E19").Value = "541 (4.9)"
$ws.Range("D20").Value = "1070 (41.9)"
$ws.Range("E20").Value = "5342 (48.1)"
$ws.Range("D21").Value = "1489 (58.3)"
$ws.Range("E21").Value = "6731 (60.6)"
$ws.Range("D22").Value = "415 (16.3)"
$ws.Range("E22").Value = "1165 (10.5)"
$ws.Range("D23").Value = "1275 (49.9)"
$ws.Range("E23").Value = "6324 (56.9)"
$ws.Range("D24").Value = "1582 (62.0)"
$ws.Range("E24").Value = "6691 (60.2)"
$ws.Range("D25").Value = "279 (10.9)"
$ws.Range("E25").Value = "1370 (12.3)"
$ws.Range("D26").Value = "2490 (97.5)"
$ws.Range("E26").Value = "10866 (97.8)"
$ws.Range("D27").Value = "1782 (69.8)"
$ws.Range("E27").Value = "7281 (65.5)"
$ws.Range("D28").Value = "955 (37.4)"
$ws.Range("E28").Value = "4064 (36.6)"
$ws.Range("D29").Value = "584 (22.9)"
$ws.Range("E29").Value = "2754 (24.8)"
$ws.Range("D30").Value = "35 (1.4)"
$ws.Range("E30").Value = "171 (1.5)"
$ws.Range("D31").Value = "784 (30.7)"
$ws.Range("E31").Value = "4313 (38.8)"
$ws.Range("D33").Value = "20 (0.8)"
$ws.Range("E33").Value = "70 (0.6)"
$ws.Range("D34").Value = "168 (6.6)"
$ws.Range("E34").Value = "549 (4.9)"
$ws.Range("D35").Value = "66 (2.6)"
$ws.Range("E35").Value = "197 (1.8)"
$ws.Range("D36").Value = "278 (10.9)"
$ws.Range("E36").Value = "490 (4.4)"
$ws.Range("D37").Value = "2018 (79.0)"
$ws.Range("E37").Value = "9805 (88.2)"
$ws.Range("D38").Value = "76 (3.0)"
$ws.Range("E38").Value = "314 (2.8)"
$ws.Range("D39").Value = "1022 (40.0)"
$ws.Range("E39").Value = "3090 (27.8)"
$ws.Range("D40").Value = "1455 (57.0)"
$ws.Range("E40").Value = "7709 (69.4)"
$ws.Range("D41").Value = "120 (4.7)"
$ws.Range("E41").Value = "499 (4.5)"
$ws.Range("D42").Value = "146 (5.7)"
$ws.Range("E42").Value = "523 (4.7)"
$ws.Range("E43").Value = "88 (0.8)"
$ws.Range("E44").Value = "23 (0.2)"
$ws.Range("E46").Value = "68 [58,79]"
$ws.Range("C47").Value = 11081
$ws.Range("D47").Value = "5.46 [3.46,9.94]"
$ws.Range("E47").Value = "5.79 [3.54,10.29]"
$ws.Range("C48").Value = 2585
$ws.Range("D48").Value = "4.23 [2.88,8.04]"
$ws.Range("E48").Value = "4.13 [2.83,7.48]"
$ws.Range("C49").Value = 11081
$ws.Range("E49").Value = "11.00 [6.00,19.00]"
$ws.Range("C50").Value = 2585
$ws.Range("D52").Value = "6 [4,9]"
$ws.Range("E52").Value = "6 [4,8]"
$ws.Range("C53").Value = 4933
$ws.Range("C54").Value = 30
$ws.Range("C55").Value = 5308
$ws.Range("C56").Value = 22
$ws.Range("C57").Value = 27
$ws.Range("C59").Value = 4942
$ws.Range("D59").Value = "748 [250,1500]"
$ws.Range("E59").Value = "681 [220,1500]"
$ws.Range("C60").Value = 310
$ws.Range("D60").Value = "2915 [1043,6446]"
$ws.Range("E60").Value = "2546 [870,5886]"
$ws.Range("C61").Value = 310
$ws.Range("D61").Value = "521.7 [240.1,984.1]"
$ws.Range("E61").Value = "487.9 [213.9,917.1]"
$ws.Range("C62").Value = 7992
$ws.Range("C63").Value = 5446
$ws.Range("D63").Value = "43.0 [18.0,97.0]"
$ws.Range("E63").Value = "35.0 [15.0,86.0]"
$ws.Range("C64").Value = 5446
$ws.Range("D64").Value = "0.31 [0.15,0.52]"
$ws.Range("E64").Value = "0.27 [0.13,0.46]"
$ws.Range("C65").Value = 5446
$ws.Range("C66").Value = 12086
$ws.Range("D66").Value = "18.0 [4.0,48.5]"
$ws.Range("E66").Value = "27.0 [6.0,69.0]"
$ws.Range("C67").Value = 5950
$ws.Range("D67").Value = "3.0 [1.0,13.0]"
$ws.Range("C68").Value = 5950
$ws.Range("D68").Value = "32.0 [11.0,72.0]"
$ws.Range("E68").Value = "31.0 [11.0,68.0]"
$ws.Range("C69").Value = 5950
$ws.Range("D69").Value = "0.24 [0.08,0.49]"
$ws.Range("E69").Value = "0.25 [0.09,0.48]"
$ws.Range("C70").Value = 28
$ws.Range("D70").Value = "19.5 [17.0,22.7]"
$ws.Range("E70").Value = "19.1 [16.8,22.1]"
$ws.Range("C71").Value = 22
$ws.Range("D71").Value = "77.2 [71.0,85.2]"
$ws.Range("E71").Value = "74.7 [69.4,81.5]"
$ws.Range("C72").Value = 720
$ws.Range("C73").Value = 24
$ws.Range("D73").Value = "97.9 [96.3,99.2]"
$ws.Range("C74").Value = 22
$ws.Range("D74").Value = "87.8 [76.6,100.8]"
$ws.Range("E74").Value = "86.0 [75.9,98.2]"
$ws.Range("C75").Value = 4140
$ws.Range("D75").Value = "88.0 [68.0,123.0]"
$ws.Range("E75").Value = "88.0 [71.0,117.0]"
$ws.Range("C76").Value = 4140
$ws.Range("C77").Value = 2339
$ws.Range("C78").Value = 69
$ws.Range("D78").Value = "155.0 [122.0,217.0]"
$ws.Range("C79").Value = 29
$ws.Range("C80").Value = 36
$ws.Range("C81").Value = 13311
$ws.Range("D81").Value = "22.2 [13.6,31.4]"
$ws.Range("E81").Value = "22.4 [12.5,35.6]"
$ws.Range("C82").Value = 1819
$ws.Range("D82").Value = "9.8 [8.3,11.4]"
$ws.Range("C83").Value = 8836
$ws.Range("D83").Value = "227.0 [153.0,361.0]"
$ws.Range("E83").Value = "230.0 [165.0,348.0]"
$ws.Range("C84").Value = 811
